$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.886.31"
$ws.Range("E2").Value = "  -3.88%  "

# Row 3
$ws.Range("D3").Value = "2.290.99"
$ws.Range("E3").Value = "  -4.93%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.03%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.01%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.37%  "

# Row 9
$ws.Range("D9").Value = "2.291.55"
$ws.Range("E9").Value = "  -4.79%  "

# Row 10
$ws.Range("E10").Value = "  -5.36%  "

# Row 11
$ws.Range("E11").Value = "  -2.73%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.148"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.51%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.332"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.79%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.39%  "

# Row 15
$ws.Range("D15").Value = "2.702.24"

# Row 16
$ws.Range("D16").Value = "57.935.26"
$ws.Range("E16").Value = "  -3.68%  "

# Row 17
$ws.Range("E17").Value = "  -4.68%  "

# Row 18
$ws.Range("D18").Value = "2.292.32"
$ws.Range("E18").Value = "  -4.15%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.67%  "

# Row 20
$ws.Range("E20").Value = "  -6.41%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "313.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.42%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.35%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.166"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.57%  "

# Row 26
$ws.Range("E26").Value = "  -0.14%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.44%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.34%  "

# Row 29
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.07%  "

# Row 30
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.29%  "

# Row 31
$ws.Range("D31").Value = "0.0₃0722"
$ws.Range("E31").Value = "  -6.61%  "

# Row 32
$ws.Range("E32").Value = "  -0.32%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.49%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.379"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.74%  "

# Row 35
$ws.Range("E35").Value = "  +0.00%  "

# Row 36
$ws.Range("E36").Value = "  -3.75%  "

# Row 37
$ws.Range("E37").Value = "  +0.02%  "

# Row 38
$ws.Range("E38").Value = "  -7.44%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.33%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.81%  "

# Row 41
$ws.Range("E41").Value = "  -6.83%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "142.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.96%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "290.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -10.17%  "

# Row 44
$ws.Range("E44").Value = "  -4.37%  "

# Row 45
$ws.Range("E45").Value = "  -2.72%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0499"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.43%  "

# Row 47
$ws.Range("E47").Value = "  -3.52%  "

# Row 48
$ws.Range("E48").Value = "  -8.16%  "

# Row 49
$ws.Range("E49").Value = "  -4.79%  "

# Row 50
$ws.Range("E50").Value = "  -3.79%  "

# Row 51
$ws.Range("E51").Value = "  -0.92%  "
